# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) for every data row (rows 2-31) was stored as the
# sheet's own file-name-derived label "5-19-2007-08"; it needs to become the
# correct ISO game date "2008-05-19".
#
# NOTE: "2008-05-19" looks like a date to Excel's text parser, so a plain
# Range.Value assignment would silently convert the cell to a date serial
# number. To keep it as literal text (matching the original inline-string
# cell), the range is temporarily formatted as Text ("@") while the value is
# written, and its original (default/"Normal") style is restored immediately
# afterwards so no visible formatting changes remain on the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateRange = $ws.Range("BF2:BF31")

# Remember the current style so we can restore it once the text is safely in.
$originalStyle = $dateRange.Style

# Force text interpretation so "2008-05-19" is stored as a literal string,
# not auto-converted into a date serial value.
$dateRange.NumberFormat = "@"
$dateRange.Value = "2008-05-19"

# Restore the original cell style/format (the date fix should not change
# how the cells look).
$dateRange.Style = $originalStyle
